# From v1.2.3 to v1.2.4
# Reorders the distinguishing step (3rd step row) content among TC3, TC4 and
# TC5 so that:
#   TC3's step now reads "Chefe Clica em excluir comprovante." / "SYSTEM Exclui o comprovante."
#   TC4's step now reads "Chefe Clica em visualizar comprovante." / "SYSTEM Exibe modal com o comprovante."
#   TC5's step now reads "Chefe Clica para detalhar a solicitação de diária." / "SYSTEM Apresenta a tela de Detalhar Diárias"
# The TC header rows (TC3/TC4/TC5 labels) and every other row stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 32 -> becomes the "excluir comprovante" step (previously under TC4, row 41)
$ws.Range("B32").Value = "Chefe Clica em excluir comprovante."
$ws.Range("D32").Value = "SYSTEM Exclui o comprovante."

# Row 41 -> becomes the "visualizar comprovante" step (previously under TC5, row 50)
$ws.Range("B41").Value = "Chefe Clica em visualizar comprovante."
$ws.Range("D41").Value = "SYSTEM Exibe modal com o comprovante."

# Row 50 -> becomes the "detalhar a solicitação de diária" step (previously under TC3, row 32)
$ws.Range("B50").Value = "Chefe Clica para detalhar a solicitação de diária."
$ws.Range("D50").Value = "SYSTEM Apresenta a tela de Detalhar Diárias"
